$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full table data (rows 2..15), matching the updated reference-sequence table.
$data = @(
    @("OQ818319*", "E. dupreanum sapovirus 1", 90, 67.180000000000007, "AQQ78883.1", 1, 79.819999999999993, "LC504397.1"),
    @("PP766459", "E. dupreanum sapovirus 1", 77, 99.9, "AQQ78883.1", 1, 79.819999999999993, "LC504397.1"),
    @("OQ818340*", "E. dupreanum sapovirus 2", 99, 67.97, "AQQ78883.1", 63, 71.400000000000006, "KX759619.1"),
    @("PP766461*", "E. dupreanum sapovirus 3", 96, 90.13, "AQQ78883.1", 92, 80.7, "KX759623.1"),
    @("PP766460*", "E. dupreanum sapovirus 4", 70, 75.260000000000005, "AQQ78883.1", "No result", "No result", "No result"),
    @("OQ818345*", "R. madagascariensis sapovirus 1", 99, 68.19, "XBH24168.1", "No result", "No result", "No result"),
    @("OQ818347*", "R. madagascariensis sapovirus 2", 99, 89.63, "XBH24156.1", 100, 77.94, "PP712001.1"),
    @("PP766470*", "R. madagascariensis sapovirus 2", 99, 90.86, "XBH24177.1", 99, 78.510000000000005, "PP712015.1"),
    @("PP766473*", "R. madagascariensis sapovirus 2", 99, 90.59, "XBH24177.1", 97, 77.28, "PP712015.1"),
    @("PP766474*", "R. madagascariensis sapovirus 2", 100, 94.55, "XBH24177.1", 68, 71.739999999999995, "PP712026.1"),
    @("PP766476*", "R. madagascariensis sapovirus 2", 99, 86.05, "XBH24177.1", 19, 73.709999999999994, "KX759623.1"),
    @("PP766477*", "R. madagascariensis sapovirus 2", 53, 81.319999999999993, "XBH24178.1", 72, 76.22, "PP712001.1"),
    @("OQ818348*", "R. madagascariensis sapovirus 3", 99, 84.12, "XBH24163.1", 94, 75, "PP712006.1"),
    @("PP766468*", "R. madagascariensis sapovirus 3", 100, 85, "XBH24156.1", 87, 75.67, "PP712033.1")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $row++
}

$ws.Range("F21").Select()
